$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44511
$ws.Range("M2").Value = 15
$ws.Range("N2").Value = 22000
$ws.Range("O2").Value = 22000
$ws.Range("P2").Value = 22000
$ws.Range("Q2").Value = "`$/caja 15 kilos granel"
$ws.Range("S2").Value = 1467
$ws.Range("T2").Value = 15
$ws.Range("D3").Value = 45085
$ws.Range("M3").Value = 110
$ws.Range("D4").Value = 45079
$ws.Range("M4").Value = 100
$ws.Range("N4").Value = 18000
$ws.Range("O4").Value = 18000
$ws.Range("P4").Value = 18000
$ws.Range("S4").Value = 1000
$ws.Range("D5").Value = 45089
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 100
$ws.Range("N5").Value = 16000
$ws.Range("O5").Value = 16000
$ws.Range("P5").Value = 16000
$ws.Range("Q5").Value = "`$/caja 18 kilos granel"
$ws.Range("S5").Value = 889
$ws.Range("T5").Value = 18
$ws.Range("D6").Value = 44392
$ws.Range("L6").Value = "Especial"
$ws.Range("M6").Value = 500
$ws.Range("N6").Value = 7000
$ws.Range("O6").Value = 7000
$ws.Range("P6").Value = 7000
$ws.Range("Q6").Value = "`$/bandeja 8 kilos"
$ws.Range("S6").Value = 875
$ws.Range("T6").Value = 8
$ws.Range("D7").Value = 44418
$ws.Range("M7").Value = 100
$ws.Range("N7").Value = 8000
$ws.Range("O7").Value = 8000
$ws.Range("P7").Value = 8000
$ws.Range("Q7").Value = "`$/caja 15 kilos granel"
$ws.Range("S7").Value = 533
$ws.Range("T7").Value = 15
$ws.Range("D8").Value = 45086
$ws.Range("M8").Value = 80
$ws.Range("N8").Value = 16000
$ws.Range("O8").Value = 16000
$ws.Range("P8").Value = 16000
$ws.Range("S8").Value = 889
$ws.Range("D10").Value = 45090
$ws.Range("M10").Value = 140
$ws.Range("N10").Value = 16000
$ws.Range("O10").Value = 16000
$ws.Range("P10").Value = 16000
$ws.Range("Q10").Value = "`$/caja 18 kilos granel"
$ws.Range("S10").Value = 889
$ws.Range("T10").Value = 18
$ws.Range("D11").Value = 44208
$ws.Range("L11").Value = "Especial"
$ws.Range("M11").Value = 70
$ws.Range("N11").Value = 24000
$ws.Range("O11").Value = 24000
$ws.Range("P11").Value = 24000
$ws.Range("Q11").Value = "`$/caja 15 kilos granel"
$ws.Range("S11").Value = 1600
$ws.Range("T11").Value = 15
$ws.Range("D12").Value = 44264
$ws.Range("L12").Value = "Calibre 100"
$ws.Range("M12").Value = 50
$ws.Range("N12").Value = 20000
$ws.Range("O12").Value = 20000
$ws.Range("P12").Value = 20000
$ws.Range("Q12").Value = "`$/caja 18 kilos embalada"
$ws.Range("S12").Value = 1111
$ws.Range("D15").Value = 44601
$ws.Range("M15").Value = 30
$ws.Range("N15").Value = 28000
$ws.Range("O15").Value = 28000
$ws.Range("P15").Value = 28000
$ws.Range("S15").Value = 1556
$ws.Range("D16").Value = 45092
$ws.Range("M16").Value = 220
$ws.Range("N16").Value = 16000
$ws.Range("P16").Value = 16000
$ws.Range("S16").Value = 889
$ws.Range("D17").Value = 44411
$ws.Range("L17").Value = "Primera"
$ws.Range("M17").Value = 210
$ws.Range("Q17").Value = "`$/bandeja 8 kilos"
$ws.Range("S17").Value = 1000
$ws.Range("T17").Value = 8
$ws.Range("D18").Value = 45093
$ws.Range("M18").Value = 170
$ws.Range("N18").Value = 15000
$ws.Range("O18").Value = 16000
$ws.Range("P18").Value = 15471
$ws.Range("Q18").Value = "`$/caja 18 kilos granel"
$ws.Range("S18").Value = 860
$ws.Range("T18").Value = 18
$ws.Range("D19").Value = 44217
$ws.Range("M19").Value = 55
$ws.Range("N19").Value = 18000
$ws.Range("O19").Value = 18000
$ws.Range("P19").Value = 18000
$ws.Range("S19").Value = 1000
$ws.Range("D20").Value = 45083
$ws.Range("L20").Value = "Primera"
$ws.Range("M20").Value = 55
$ws.Range("N20").Value = 16000
$ws.Range("O20").Value = 16000
$ws.Range("P20").Value = 16000
$ws.Range("Q20").Value = "`$/caja 18 kilos granel"
$ws.Range("S20").Value = 889
